$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Fgf13"
$ws.Range("C2").Value = "Scn8a"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03538133333333333
$ws.Range("H2").Value = 0.106144
$ws.Range("I2").Value = 0.0113454327482436
$ws.Range("J2").Value = 0.01134543274824361
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8653686666666666
$ws.Range("N2").Value = 2.596106
$ws.Range("O2").Value = 0.7337618712085689
$ws.Range("P2").Value = 0.7337618712085689
$ws.Range("Q2").Value = 0.03061789725155555
$ws.Range("R2").Value = 0.275561075264
$ws.Range("S2").Value = 0.008324845963022203
$ws.Range("T2").Value = 0.008324845963022205

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Fgf13"
$ws.Range("C3").Value = "Scn8a"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03538133333333333
$ws.Range("H3").Value = 0.106144
$ws.Range("I3").Value = 0.0113454327482436
$ws.Range("J3").Value = 0.01134543274824361
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.3139903333333333
$ws.Range("N3").Value = 0.941971
$ws.Range("O3").Value = 0.2662381287914311
$ws.Range("P3").Value = 0.2662381287914311
$ws.Range("Q3").Value = 0.01110939664711111
$ws.Range("R3").Value = 0.099984569824
$ws.Range("S3").Value = 0.0030205867852214
$ws.Range("T3").Value = 0.003020586785221401

$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Fgf13"
$ws.Range("C4").Value = "Scn8a"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.083171666666666
$ws.Range("H4").Value = 9.249514999999999
$ws.Range("I4").Value = 0.9886545672517564
$ws.Range("J4").Value = 0.9886545672517564
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8653686666666666
$ws.Range("N4").Value = 2.596106
$ws.Range("O4").Value = 0.7337618712085689
$ws.Range("P4").Value = 0.7337618712085689
$ws.Range("Q4").Value = 2.668080154287777
$ws.Range("R4").Value = 24.01272138858999
$ws.Range("S4").Value = 0.7254370252455468
$ws.Range("T4").Value = 0.7254370252455468

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Fgf13"
$ws.Range("C5").Value = "Scn8a"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.083171666666666
$ws.Range("H5").Value = 9.249514999999999
$ws.Range("I5").Value = 0.9886545672517564
$ws.Range("J5").Value = 0.9886545672517564
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3139903333333333
$ws.Range("N5").Value = 0.941971
$ws.Range("O5").Value = 0.2662381287914311
$ws.Range("P5").Value = 0.2662381287914311
$ws.Range("Q5").Value = 0.9680860993405555
$ws.Range("R5").Value = 8.712774894064999
$ws.Range("S5").Value = 0.2632175420062097
$ws.Range("T5").Value = 0.2632175420062097
